$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-25 17:18:39'
$ws.Range('E3').Value = '2026-02-25 17:18:42'
$ws.Range('K3').Value = '16.6 MJ/m2'
$ws.Range('O3').Value = '3.8 °C'
$ws.Range('E4').Value = '2026-02-25 17:18:44'
$ws.Range('J4').Value = '1021.5 hPa'
$ws.Range('K4').Value = '5.1 MJ/m2'
$ws.Range('O4').Value = '7.9 °C'
$ws.Range('E5').Value = '2026-02-25 17:18:47'
$ws.Range('H5').Value = '''24%'
$ws.Range('K5').Value = '15.7 MJ/m2'
$ws.Range('E6').Value = '2026-02-25 17:18:49'
$ws.Range('J6').Value = '1021.3 hPa'
$ws.Range('O6').Value = '10.7 °C'
$ws.Range('E7').Value = '2026-02-25 17:18:52'
$ws.Range('H7').Value = '''88%'
$ws.Range('K7').Value = '11.2 MJ/m2'
$ws.Range('O7').Value = '12.5 °C'
$ws.Range('E8').Value = '2026-02-25 17:18:55'
$ws.Range('J8').Value = '1020.4 hPa'
$ws.Range('K8').Value = '10.9 MJ/m2'
$ws.Range('O8').Value = '12.3 °C'
$ws.Range('E9').Value = '2026-02-25 17:18:57'
$ws.Range('O9').Value = '9.8 °C'
$ws.Range('E10').Value = '2026-02-25 17:19:00'
$ws.Range('K10').Value = '8.3 MJ/m2'
$ws.Range('O10').Value = '9.5 °C'
$ws.Range('E11').Value = '2026-02-25 17:19:03'
$ws.Range('O11').Value = '9.1 °C'
$ws.Range('E12').Value = '2026-02-25 17:19:05'
$ws.Range('O12').Value = '9.6 °C'
$ws.Range('E13').Value = '2026-02-25 17:19:08'
$ws.Range('H13').Value = '''67%'
$ws.Range('J13').Value = '1022.9 hPa'
$ws.Range('K13').Value = '15.1 MJ/m2'
$ws.Range('O13').Value = '6.1 °C'
$ws.Range('E14').Value = '2026-02-25 17:19:11'
$ws.Range('K14').Value = '11.0 MJ/m2'
$ws.Range('O14').Value = '10.2 °C'
$ws.Range('E15').Value = '2026-02-25 17:19:13'
$ws.Range('O15').Value = '9.6 °C'
$ws.Range('E16').Value = '2026-02-25 17:19:16'
$ws.Range('H16').Value = '''27%'
$ws.Range('N16').Value = '2.5 °C 16:49 TU'
$ws.Range('O16').Value = '3.4 °C'
$ws.Range('E17').Value = '2026-02-25 17:19:18'
$ws.Range('K17').Value = '17.9 MJ/m2'
$ws.Range('O17').Value = '9.8 °C'
$ws.Range('E18').Value = '2026-02-25 17:19:21'
$ws.Range('K18').Value = '6.2 MJ/m2'
$ws.Range('O18').Value = '9.8 °C'
$ws.Range('E19').Value = '2026-02-25 17:19:24'
$ws.Range('K19').Value = '15.9 MJ/m2'
$ws.Range('E20').Value = '2026-02-25 17:19:26'
$ws.Range('H20').Value = '''50%'
$ws.Range('K20').Value = '16.6 MJ/m2'
$ws.Range('E21').Value = '2026-02-25 17:19:28'
$ws.Range('J21').Value = '1021.5 hPa'
$ws.Range('K21').Value = '15.7 MJ/m2'
$ws.Range('O21').Value = '9.5 °C'
$ws.Range('E22').Value = '2026-02-25 17:19:31'
$ws.Range('K22').Value = '16.5 MJ/m2'
$ws.Range('E23').Value = '2026-02-25 17:19:33'
$ws.Range('H23').Value = '''30%'
$ws.Range('K23').Value = '16.4 MJ/m2'
$ws.Range('E24').Value = '2026-02-25 17:19:36'
$ws.Range('H24').Value = '''77%'
$ws.Range('K24').Value = '15.2 MJ/m2'
$ws.Range('L24').Value = '17.3 km/h - 199º 16:40 TU'
$ws.Range('O24').Value = '10.6 °C'
$ws.Range('E25').Value = '2026-02-25 17:19:39'
$ws.Range('K25').Value = '17.3 MJ/m2'
$ws.Range('E26').Value = '2026-02-25 17:19:41'
$ws.Range('H26').Value = '''45%'
$ws.Range('K26').Value = '16.5 MJ/m2'
$ws.Range('O26').Value = '10.8 °C'
$ws.Range('E27').Value = '2026-02-25 17:19:44'
$ws.Range('K27').Value = '16.5 MJ/m2'
$ws.Range('L27').Value = '24.5 km/h - 183º 16:56 TU'
$ws.Range('E28').Value = '2026-02-25 17:19:47'
$ws.Range('J28').Value = '1021.4 hPa'
$ws.Range('K28').Value = '13.2 MJ/m2'
$ws.Range('O28').Value = '8.4 °C'
$ws.Range('E29').Value = '2026-02-25 17:19:49'
$ws.Range('K29').Value = '6.6 MJ/m2'
$ws.Range('E30').Value = '2026-02-25 17:19:52'
$ws.Range('O30').Value = '10.3 °C'
$ws.Range('E31').Value = '2026-02-25 17:19:55'
$ws.Range('J31').Value = '1021.0 hPa'
$ws.Range('E32').Value = '2026-02-25 17:19:57'
$ws.Range('K32').Value = '16.5 MJ/m2'
$ws.Range('O32').Value = '9.9 °C'
$ws.Range('E33').Value = '2026-02-25 17:20:00'
$ws.Range('H33').Value = '''51%'
$ws.Range('J33').Value = '1021.1 hPa'
$ws.Range('K33').Value = '15.6 MJ/m2'
$ws.Range('O33').Value = '8.2 °C'
$ws.Range('E34').Value = '2026-02-25 17:20:03'
$ws.Range('O34').Value = '4.0 °C'
$ws.Range('E35').Value = '2026-02-25 17:20:05'
$ws.Range('K35').Value = '16.6 MJ/m2'
$ws.Range('E36').Value = '2026-02-25 17:20:08'
$ws.Range('K36').Value = '4.4 MJ/m2'
$ws.Range('O36').Value = '10.9 °C'
$ws.Range('E37').Value = '2026-02-25 17:20:11'
$ws.Range('O37').Value = '6.9 °C'
$ws.Range('E38').Value = '2026-02-25 17:20:13'
$ws.Range('H38').Value = '''88%'
$ws.Range('K38').Value = '6.4 MJ/m2'
$ws.Range('O38').Value = '9.0 °C'
$ws.Range('E39').Value = '2026-02-25 17:20:16'
$ws.Range('E40').Value = '2026-02-25 17:20:18'
$ws.Range('J40').Value = '1021.6 hPa'
$ws.Range('O40').Value = '9.7 °C'
$ws.Range('E41').Value = '2026-02-25 17:20:21'
$ws.Range('J41').Value = '1020.5 hPa'
$ws.Range('K41').Value = '9.2 MJ/m2'
$ws.Range('E42').Value = '2026-02-25 17:20:24'
$ws.Range('O42').Value = '11.3 °C'
$ws.Range('E43').Value = '2026-02-25 17:20:26'
$ws.Range('H43').Value = '''70%'
$ws.Range('K43').Value = '15.9 MJ/m2'
$ws.Range('O43').Value = '9.6 °C'
$ws.Range('E44').Value = '2026-02-25 17:20:29'
$ws.Range('K44').Value = '15.9 MJ/m2'
$ws.Range('O44').Value = '2.1 °C'
$ws.Range('E45').Value = '2026-02-25 17:20:31'
$ws.Range('J45').Value = '1019.4 hPa'
$ws.Range('O45').Value = '11.2 °C'
$ws.Range('E46').Value = '2026-02-25 17:20:34'
$ws.Range('H46').Value = '''82%'
$ws.Range('K46').Value = '15.1 MJ/m2'
$ws.Range('O46').Value = '9.2 °C'
